{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Locate the two paragraphs we need to touch:\n//  - targetPara: the body paragraph ending in\n//    \"Use Heading 2 and 3 and paragraph to conversate about the issue of the Heading 1\"\n//  - headingPara: the Heading 1 styled paragraph whose text is exactly \"Heading 1\"\nlet targetPara = null;\nlet headingPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Use Heading 2 and 3 and paragraph\") !== -1) {\n    targetPara = p;\n  }\n  if (p.style === \"Heading 1\" && p.text === \"Heading 1\") {\n    headingPara = p;\n  }\n}\n\nif (!targetPara || !headingPara) {\n  throw new Error(\"Could not locate expected paragraphs to edit.\");\n}\n\n// Helper: search within a scoped range, insert `insertStr` right after the\n// LAST match of `searchText` (there can be earlier, unrelated matches of the\n// same substring elsewhere in the paragraph).\nasync function insertAfterLastMatch(scopeRange, searchText, insertStr) {\n  const results = scopeRange.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[results.items.length - 1];\n  found.insertText(insertStr, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// --- Change 1 -------------------------------------------------------------\n// \"Use Heading 2 and 3 and paragraph to conversate about the issue of the Heading 1\"\n// becomes\n// \"Use Heading type 2 and type 3 and paragraph to conversate about the issue of the Heading type 1\"\nawait insertAfterLastMatch(targetPara, \"Use Heading \", \"type \");\nawait insertAfterLastMatch(targetPara, \"2 and \", \"type \");\nawait insertAfterLastMatch(targetPara, \"the issue of the Heading \", \"type \");\n\n// --- Change 2 -------------------------------------------------------------\n// The \"Heading 1\" paragraph text is replaced with a pointer to the reference\n// document, and a new empty Heading 1 paragraph is appended after it.\nheadingPara.insertText(\n  \"Please read \\u201CAssignment Reference Documents.pdf\\u201D in the HUMAN RESOURCES folder of this git repo.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst endOfHeadingPara = headingPara.getRange(Word.RangeLocation.end);\nendOfHeadingPara.insertText(\"\\n\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# Curly quotes used in the replacement text below.\n$lq = [char]0x201C   # \u201c\n$rq = [char]0x201D   # \u201d\n\n# Locate the two paragraphs we need to touch:\n#  - $targetPara: the body paragraph ending in\n#    \"Use Heading 2 and 3 and paragraph to conversate about the issue of the Heading 1\"\n#  - $headingPara: the Heading 1 styled paragraph whose text is exactly \"Heading 1\"\n$targetPara = $null\n$headingPara = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    $styleName = $p.Range.Style.NameLocal\n    if ($text -like \"*Use Heading 2 and 3 and paragraph*\") {\n        $targetPara = $p\n    }\n    if ($styleName -eq \"Heading 1\" -and $text.TrimEnd([char]13) -eq \"Heading 1\") {\n        $headingPara = $p\n    }\n}\n\n# --- Change 1 ---------------------------------------------------------\n# \"Use Heading 2 and 3 and paragraph to conversate about the issue of the Heading 1\"\n# becomes\n# \"Use Heading type 2 and type 3 and paragraph to conversate about the issue of the Heading type 1\"\n$oldSentence = \"Use Heading 2 and 3 and paragraph to conversate about the issue of the Heading 1\"\n$newSentence = \"Use Heading type 2 and type 3 and paragraph to conversate about the issue of the Heading type 1\"\n\n$scope1 = $targetPara.Range.Duplicate\n$scope1.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 1) | Out-Null\n\n# --- Change 2 ---------------------------------------------------------\n# The \"Heading 1\" paragraph text is replaced with a pointer to the reference\n# document, and a new empty Heading 1 paragraph is appended after it.\n$replacement = \"Please read \" + $lq + \"Assignment Reference Documents.pdf\" + $rq + \" in the HUMAN RESOURCES folder of this git repo.\"\n\n$scope2 = $headingPara.Range.Duplicate\n$scope2.Find.Execute(\"Heading 1\", $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 1) | Out-Null\n\n$headingPara.Range.InsertParagraphAfter()\n"}
